$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "9 MESES DE TRAINING" block (rows 20-22) so it references its own
# experiment outputs (TS9311 / HT9411 / ZZ9411 / ZZ9412) instead of the ones
# copied over from the "3 MESES DE TRAINING" block (TS9310 / HT9410 / ZZ9410 / ZZ9411).

$ws.Range("E20").Value = "exp/TS9311/dataset_training.csv.gz"
$ws.Range("F20").Value = "exp/HT9411/dataset_training.csv.gz"

$ws.Range("E21").Value = "exp/HT9411/dataset_training.csv.gz"
$ws.Range("F21").Value = "exp/ZZ9411"

$ws.Range("E22").Value = "exp/HT9411/dataset_training.csv.gz"
$ws.Range("F22").Value = "exp/ZZ9412"

# Move the active selection to A19, matching the state at save time.
$ws.Range("A19").Select()

$wb.Save()
